$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.340.11'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '1.895.22'
$ws.Range('E3').Value = '  +1.60%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.51'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('E6').Value = '  +3.23%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.86'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('E9').Value = '  +5.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.35'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +11.02%  '
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.93'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.792'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +12.14%  '
$ws.Range('D15').Value = '2.171.20'
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('E16').Value = '  +4.58%  '
$ws.Range('D17').Value = '1.919.14'
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').Value = '35.419.66'
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').Value = '0.0₃0828'
$ws.Range('E20').Value = '  +2.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.42'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.94'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.22'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +7.54%  '
$ws.Range('E24').Value = '  +8.26%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.51'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.22%  '
$ws.Range('E28').Value = '  +3.08%  '
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  +6.71%  '
$ws.Range('E32').Value = '  +5.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.22'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('E34').Value = '  +25.26%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  -13.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.853'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.68%  '
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0718'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +8.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0225'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +8.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.93'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.95'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('B44').Value = 'Gas'
$ws.Range('C44').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.66'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +14.93%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.337.59'
$ws.Range('E45').Value = '  +4.85%  '
$ws.Range('E46').Value = '  +3.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0811'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.60'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.91%  '
